$d = $word.ActiveDocument

# Locate the paragraph that ends with "外交手段" (Tact / 外交手段) so we can
# insert the new "As far as / 據我所知" paragraph right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*外交手段*") {
        $target = $p
    }
}

# Insert a new paragraph right after the target paragraph.
$target.Range.InsertParagraphAfter()

# The newly created paragraph is the one right after $target.
$newPara = $target.Next()
$newRange = $newPara.Range

# First run: "As far as " (western text, no special font hint)
$newRange.InsertAfter("As far as ")

# Second run: "據我所知" with an eastAsia font hint, placed right after the
# first run and before the paragraph mark.
$insertPoint = $d.Range($newRange.End - 1, $newRange.End - 1)
$insertPoint.InsertAfter("據我所知")
$insertPoint.Font.NameFarEast = "DengXian"

# --- styles.xml: add the three missing latent style exceptions ---
$d.Styles.AutoUpdate = $d.Styles.AutoUpdate
